# Actualización automática hashcode dom jun 30 02:00:40 CEST 2019
# Updates the hashcode values (column B) for a set of rows identified by
# their ID (column A) in the "hashcode.csv" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @{
    100 = "85819c9b0ee864700a6fb3abf7b62758"  # 04-040021TM
    104 = "afc45b0ea45fcd2114d8102997488408"  # 04-040021TP
    113 = "956b266fd844e9f3fca2194ee278fadb"  # 04-040021TC
    122 = "d15ca3c8fb72fbbd22db7c2394f28a69"  # 04-040014TC
    164 = "0a80cf60deec27272e68c8141fbee685"  # 04-040021A
    230 = "a7ccd9496d18261177551264266f67e7"  # 04-040014TP
    233 = "380c5e4c6ed05e85df43317f9a0cfa66"  # 04-040014TM
    331 = "d9986ed4380897b50d61c0803314de7c"  # 04-040018TP
    342 = "052d5b4453144717d9154004c40aed09"  # 04-040018TC
    343 = "9c8e173b79f48d63f00af95644862e76"  # 04-040018TM
    381 = "426758b07b194188b97fe09b886f440d"  # 01-010073A
    419 = "0841f66eec1f7caf51680bed6f5054c6"  # 05-0709-070905BTC
    458 = "62f05aaa5756711c583f9c74bdffd409"  # 01-010073TP
    477 = "e1b8840a7130774ea1c4a2335241f85b"  # 01-010073TC
    619 = "bd09cfb4e9f5a5a1edc58ee2f6cbef23"  # 04-040015TC
    623 = "5df9e1ffb7ca51b90d6720532ccfee6f"  # 04-040015TP
    628 = "ae8a27b09551a4de674da30e82a0e23c"  # 04-040015TM
    779 = "babf3fd530aff2ea45435a4292853ff1"  # 04-040018A
    818 = "4c2ed9e49577e877cba8646fab52dc00"  # 04-040015A
    831 = "3ebef27ff7385eb5bb0c6c1d9dc07834"  # 04-040014A
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
